$d = $word.ActiveDocument

# 1. Collapse the ">>> your stuff after this line >>>" paragraph (currently split
#    across three runs with proofErr marks around "> your") into a single run.
#    A Find/Replace over the whole matched phrase merges the runs and drops the
#    proofErr elements, exactly like Word itself does on a literal replace.
$d.Content.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2)

# 2. Add a new paragraph after it containing the "Version management..." sentence.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfDoc = $lastPara.Range
$endOfDoc.Collapse(0)
$endOfDoc.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.InsertAfter("Version management systems are a daily reality for the software development professional.")

# Re-fetch the paragraph/range for the inserted text and apply its run formatting
# (Calibri for ascii/hAnsi/cs, English (US) language) without touching the
# paragraph mark itself.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.MoveEnd(1, -1)
$newRange.Font.Name = "Calibri"
$newRange.Font.NameBi = "Calibri"
$newRange.LanguageID = "en-US"
